$d = $word.ActiveDocument
$table = $d.Tables.Item(1)
$values = @(
    "70-10=",
    "32+52=",
    "10+49=",
    "58+11=",
    "9-9=",
    "55-43=",
    "65-57=",
    "20+0=",
    "57-21=",
    "56+16=",
    "44+2=",
    "55+31=",
    "13+38=",
    "44+47=",
    "42+52=",
    "1+86=",
    "75-20=",
    "8+67=",
    "15+4=",
    "99-11=",
    "27-5=",
    "16+32=",
    "29-20=",
    "65-21=",
    "35-24=",
    "87-31=",
    "37-32=",
    "14+15=",
    "52-45=",
    "59+16=",
    "55-18=",
    "73-21=",
    "41+28=",
    "33+44=",
    "30+9=",
    "42-41=",
    "62+14=",
    "56-24=",
    "74-50=",
    "46-16=",
    "77-74=",
    "29-10=",
    "40+28=",
    "9+15=",
    "99-28=",
    "5+67=",
    "87-59=",
    "34+25=",
    "92-84=",
    "22+24=",
    "15+9=",
    "48+10=",
    "38+38=",
    "44-17=",
    "91-37=",
    "57-10=",
    "83-58=",
    "57+3=",
    "20+22=",
    "73-28=",
    "89-22=",
    "57-0=",
    "35-6=",
    "29-22=",
    "71-19=",
    "19+48=",
    "66+12=",
    "31+0=",
    "61-8=",
    "81-54=",
    "40+1=",
    "11+80=",
    "48-20=",
    "81-38=",
    "19+76=",
    "38+16=",
    "74+21=",
    "1+69=",
    "27+28=",
    "23+67=",
    "71-42=",
    "91+0=",
    "25-8=",
    "39+23=",
    "43+17=",
    "90-80=",
    "58+0=",
    "20+63=",
    "13+12=",
    "40+42=",
    "9-0=",
    "54+36=",
    "79+16=",
    "46+35=",
    "70-8=",
    "99-81=",
    "54-38=",
    "15+48=",
    "28-3=",
    "50+34="
)

$idx = 0
$numRows = $table.Rows.Count
$numCols = $table.Columns.Count
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $table.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $values[$idx]
        $idx++
    }
}
